$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.623.23'
$ws.Range("E2").Value = '  -2.30%  '

$ws.Range("D3").Value = '1.759.38'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.02%  '

$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4329'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3607'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.45%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07563'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.111'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.073'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.51%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.231'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.93%  '

$ws.Range("D16").Value = '1.758.51'
$ws.Range("E16").Value = '  -3.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001069'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06417'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.24%  '

$ws.Range("E20").Value = '  -0.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.853'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.54%  '

$ws.Range("D23").Value = '27.671.10'
$ws.Range("E23").Value = '  -2.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.70%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.096'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.92%  '

$ws.Range("D28").Value = '1.958.11'
$ws.Range("E28").Value = '  -3.98%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.148'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.27%  '

$ws.Range("E30").Value = '  -2.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.102'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.593'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.662'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08970'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.78%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02307'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.83%  '

$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6429'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.04%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2106'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06016'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.955'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.187'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.400'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.901'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5939'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.718'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.988'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.169'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.00%  '

$ws.Range("E51").Value = '  -1.79%  '
